# Generate Report for Handback
#
# The handback finished: both language files (zh-cn / de-de) for each of the
# two tracked sources are now in sync with en-US. This fills in the
# "Latest Target File" (F) / "Latest Handback File" (G) columns on each
# language sheet (mirroring the already-known Handoff File/md name), stamps
# the "Latest Handback DateTime" (H) column with the real handback time
# (replacing the 0001-01-01 sentinel), and flips every "Status" cell from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR for RGB FF6495ED, matching the workbook's existing HyperLink style
$hyperlinkUnderline = 2      # xlUnderlineStyleSingle

function Get-HyperlinkTarget {
    param($Worksheet, $CellAddr)
    $links = @($Worksheet.Hyperlinks)
    foreach ($link in $links) {
        if ($link.Range.Address() -eq $CellAddr) {
            return $link.Address
        }
    }
    return $null
}

function Copy-CellAsHyperlink {
    param($Worksheet, $SourceAddr, $DestAddr)
    $source = $Worksheet.Range($SourceAddr)
    $dest = $Worksheet.Range($DestAddr)
    $display = $source.Value2
    $target = Get-HyperlinkTarget $Worksheet $SourceAddr
    $dest.Value = $display
    $Worksheet.Hyperlinks.Add($dest, $target, "", "", $display) | Out-Null
    $dest.Font.Underline = $hyperlinkUnderline
    $dest.Font.Color = $hyperlinkColor
}

function Update-LanguageSheet {
    param($SheetName, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 2: first tracked source file
    $ws.Range("C2").Value = $statusText
    Copy-CellAsHyperlink $ws '$A$2' "F2"
    Copy-CellAsHyperlink $ws '$D$2' "G2"
    $ws.Range("H2").Value = $HandbackDateTime

    # Row 3: second tracked source file
    $ws.Range("C3").Value = $statusText
    Copy-CellAsHyperlink $ws '$A$3' "F3"
    Copy-CellAsHyperlink $ws '$D$3' "G3"
    $ws.Range("H3").Value = $HandbackDateTime
}

Update-LanguageSheet "zh-cn" "2016-03-31 07:08:25"
Update-LanguageSheet "de-de" "2016-03-31 07:08:41"

# The Overview sheet's per-language Status cells mirror the same text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText
